$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the rate summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.73 = 6495.51 pesos`n✅ 6495.51 pesos = 1.72 = 902.04 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $text

# --- Sheet "tasas": update computed rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 577
$ws2.Range("O10").Value = 3747.91
$ws2.Range("N12").Value = 3773.99
$ws2.Range("O12").Value = 524.1
